# Apply updated "想去人数" (F column) values as scraped at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5632
$ws1.Range("F9").Value = 533

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5632
$ws4.Range("F11").Value = 533
